$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# Row 53
$ws.Range("H53").Value = 640.9
$ws.Range("I53").Value = 825.7143
$ws.Range("K53").Value = 825.7143
$ws.Range("M53").Value = -188.7143

# Row 74
$ws.Range("H74").Value = 4170.25
$ws.Range("I74").Value = 4202.4443
$ws.Range("J74").Value = 4128.857
$ws.Range("K74").Value = 4202.4443
$ws.Range("L74").Value = 4128.857
$ws.Range("M74").Value = -3266.4443
$ws.Range("N74").Value = -6000.857

# Row 77
$ws.Range("H77").Value = 4170.25
$ws.Range("I77").Value = 4202.4443
$ws.Range("J77").Value = 4128.857
$ws.Range("K77").Value = 21012.2215
$ws.Range("L77").Value = 20644.285
$ws.Range("M77").Value = -16332.2215
$ws.Range("N77").Value = -30004.285

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# Row 129
$ws.Range("H129").Value = 9604.08
$ws.Range("I129").Value = 446.16666
$ws.Range("K129").Value = 1338.49998
$ws.Range("M129").Value = 3661.50002

# Row 132
$ws.Range("H132").Value = 4193.339
$ws.Range("I132").Value = 4282.592
$ws.Range("J132").Value = 3856.923
$ws.Range("K132").Value = 12847.776
$ws.Range("L132").Value = 11570.769
$ws.Range("M132").Value = -10317.776
$ws.Range("N132").Value = -16630.769

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13016.88
$ws.Range("I32").Value = 7325.831
$ws.Range("J32").Value = 32069.521
$ws.Range("K32").Value = 7325.831
$ws.Range("L32").Value = 32069.521
$ws.Range("M32").Value = -7038.831
$ws.Range("N32").Value = -32643.521

# Row 97
$ws.Range("H97").Value = 1445.3334
$ws.Range("I97").Value = 1661.8572
$ws.Range("J97").Value = 687.5
$ws.Range("K97").Value = 1661.8572
$ws.Range("L97").Value = 687.5
$ws.Range("M97").Value = -1165.8572
$ws.Range("N97").Value = -1679.5

# Row 105
$ws.Range("H105").Value = 41999
$ws.Range("J105").Value = 41999
$ws.Range("L105").Value = 41999
$ws.Range("N105").Value = -48987

# Row 122
$ws.Range("H122").Value = 2791.7368
$ws.Range("I122").Value = 2380.4375
$ws.Range("J122").Value = 4985.3335
$ws.Range("K122").Value = 7141.3125
$ws.Range("L122").Value = 14956.0005
$ws.Range("M122").Value = -4691.3125
$ws.Range("N122").Value = -19856.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 11767453
$ws.Range("I86").Value = 15387315
$ws.Range("J86").Value = 2901.75
$ws.Range("K86").Value = 15387315
$ws.Range("L86").Value = 2901.75
$ws.Range("M86").Value = -15386192
$ws.Range("N86").Value = -5147.75

# Row 89
$ws.Range("H89").Value = 11767453
$ws.Range("I89").Value = 15387315
$ws.Range("J89").Value = 2901.75
$ws.Range("K89").Value = 76936575
$ws.Range("L89").Value = 14508.75
$ws.Range("M89").Value = -76930959
$ws.Range("N89").Value = -25740.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2089.848
$ws.Range("I31").Value = 1136.1111
$ws.Range("K31").Value = 1136.1111
$ws.Range("M31").Value = -841.1111000000001

# Row 34
$ws.Range("H34").Value = 2089.848
$ws.Range("I34").Value = 1136.1111
$ws.Range("K34").Value = 1136.1111
$ws.Range("M34").Value = -934.1111000000001

# Row 58
$ws.Range("H58").Value = 2900.3667
$ws.Range("I58").Value = 889.5625
$ws.Range("J58").Value = 5198.4287
$ws.Range("K58").Value = 889.5625
$ws.Range("L58").Value = 5198.4287
$ws.Range("M58").Value = -686.5625
$ws.Range("N58").Value = -5604.4287

# Row 122
$ws.Range("H122").Value = 1302.7
$ws.Range("I122").Value = 887.4286
$ws.Range("J122").Value = 1526.3077
$ws.Range("K122").Value = 2662.2858
$ws.Range("L122").Value = 4578.9231
$ws.Range("M122").Value = -212.2857999999997
$ws.Range("N122").Value = -9478.9231

# Row 134
$ws.Range("H134").Value = 8801.6875
$ws.Range("I134").Value = 9370.462
$ws.Range("J134").Value = 6337
$ws.Range("K134").Value = 28111.386
$ws.Range("L134").Value = 19011
$ws.Range("M134").Value = -25576.386
$ws.Range("N134").Value = -24081

# Row 136
$ws.Range("H136").Value = 2900.3667
$ws.Range("I136").Value = 889.5625
$ws.Range("J136").Value = 5198.4287
$ws.Range("K136").Value = 2668.6875
$ws.Range("L136").Value = 15595.2861
$ws.Range("M136").Value = -118.6875
$ws.Range("N136").Value = -20695.2861

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 779.01697
$ws.Range("I131").Value = 482.31818
$ws.Range("J131").Value = 955.43243
$ws.Range("K131").Value = 1446.95454
$ws.Range("L131").Value = 2866.29729
$ws.Range("M131").Value = 3593.04546
$ws.Range("N131").Value = -12946.29729

# Row 132
$ws.Range("H132").Value = 1300.625
$ws.Range("I132").Value = 506.44446
$ws.Range("J132").Value = 1777.1333
$ws.Range("K132").Value = 4558.00014
$ws.Range("L132").Value = 15994.1997
$ws.Range("M132").Value = -2028.00014
$ws.Range("N132").Value = -21054.1997

# Row 133
$ws.Range("H133").Value = 1551.625
$ws.Range("I133").Value = 1875.1428
$ws.Range("K133").Value = 5625.428400000001
$ws.Range("M133").Value = -565.4284000000007

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2424.4167
$ws.Range("I122").Value = 2621.4443
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 7864.3329
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -5414.3329
$ws.Range("N122").Value = -10400.0002

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1542.3334
$ws.Range("I46").Value = 2396.2
$ws.Range("J46").Value = 475
$ws.Range("K46").Value = 2396.2
$ws.Range("L46").Value = 475
$ws.Range("M46").Value = -2208.2
$ws.Range("N46").Value = -851

# Row 61
$ws.Range("H61").Value = 1671.579
$ws.Range("I61").Value = 1388.3334
$ws.Range("J61").Value = 2157.1428
$ws.Range("K61").Value = 1388.3334
$ws.Range("L61").Value = 2157.1428
$ws.Range("M61").Value = -1186.3334
$ws.Range("N61").Value = -2561.1428

# Row 93
$ws.Range("H93").Value = 1349.3043
$ws.Range("I93").Value = 1417.5385
$ws.Range("J93").Value = 1260.6
$ws.Range("K93").Value = 1417.5385
$ws.Range("L93").Value = 1260.6
$ws.Range("M93").Value = -169.5385000000001
$ws.Range("N93").Value = -3756.6

# Row 113
$ws.Range("H113").Value = 1671.579
$ws.Range("I113").Value = 1388.3334
$ws.Range("J113").Value = 2157.1428
$ws.Range("K113").Value = 1388.3334
$ws.Range("L113").Value = 2157.1428
$ws.Range("M113").Value = 781.6666
$ws.Range("N113").Value = -6497.1428

# Row 122
$ws.Range("H122").Value = 5459.8
$ws.Range("I122").Value = 5514.0713
$ws.Range("J122").Value = 5333.1665
$ws.Range("K122").Value = 16542.2139
$ws.Range("L122").Value = 15999.4995
$ws.Range("M122").Value = -14092.2139
$ws.Range("N122").Value = -20899.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 30375
$ws.Range("J80").Value = 30375
$ws.Range("L80").Value = 30375
$ws.Range("N80").Value = -32371

# Row 83
$ws.Range("H83").Value = 30375
$ws.Range("J83").Value = 30375
$ws.Range("L83").Value = 91125
$ws.Range("N83").Value = -101109

# Row 104
$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -26988

# Row 122
$ws.Range("H122").Value = 2450
$ws.Range("I122").Value = 2450
$ws.Range("K122").Value = 7350
$ws.Range("M122").Value = -4900
